# Add a "Prix" (price) column in B and a "Lieu" (place) column in C next to
# the existing "Noms des applications" column in A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "Prix"
$ws.Range("C1").Value = "Lieu"

# Data rows
$ws.Range("B2").Value = 19745
$ws.Range("C2").Value = "Lion"

$ws.Range("B3").Value = 1475621
$ws.Range("C3").Value = "Lile"

$ws.Range("B4").Value = 155

$ws.Range("B5").Value = 834548
